# Utilidades Mineria de Datos
# Adds two new worksheets ("001" and "002") containing Spanish/English
# country reference tables, and updates view state (active sheet, zoom,
# selection) to match the authored workbook.

$wb = $excel.ActiveWorkbook
$ws0 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Add the two new sheets, in order, right after "000"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Add($null, $ws0)
$ws1.Name = "001"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "002"

# ---------------------------------------------------------------------
# 2. Populate sheet "001" -- Estado / Año table
# ---------------------------------------------------------------------
$sheet001 = @(
    @('Estado', 'Año'),
    @('Austria', 2002),
    @('Bélgica', 2002),
    @('Bulgaria', 2002),
    @('Comunidad Helvética', 2002),
    @('Chipre', 2002),
    @('Alemania', 2002),
    @('Dinamarca', 2002),
    @('Estonia', 2002),
    @('España', 2002),
    @('Finlandia', 2002),
    @('Francia', 2002),
    @('Gran Bretaña', 2002),
    @('Hungría', 2002),
    @('Irlanda', 2002),
    @('Países Bajos', 2002),
    @('Noruega', 2002),
    @('Polonia', 2002),
    @('Portugal', 2002),
    @('Rusia', 2002),
    @('Suecia', 2002),
    @('Eslovenia', 2002),
    @('Eslovakia', 2002),
    @('Ucrania', 2002)
)

for ($i = 0; $i -lt $sheet001.Length; $i++) {
    $r = $i + 1
    $row = $sheet001[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
}

# ---------------------------------------------------------------------
# 3. Populate sheet "002" -- Abbreviation / Countries table
# ---------------------------------------------------------------------
$sheet002 = @(
    @('Abbreviation', 'Countries'),
    @('NL', 'Netherlands'),
    @('AT', 'Austria'),
    @('SI', 'Slovenia'),
    @('BG', 'Bulgary'),
    @('CZ', 'Czech Republic'),
    @('CY', 'Cyprus'),
    @('DE', 'Germany'),
    @('DK', 'Denmark'),
    @('EE', 'Estonia'),
    @('ES', 'Spain'),
    @('FI', 'Finland'),
    @('FR', 'France'),
    @('GB', 'Great Britain'),
    @('HU', 'Hungary'),
    @('IE', 'Ireland'),
    @('RO', 'Romania'),
    @('NO', 'Norway'),
    @('PL', 'Polond'),
    @('PT', 'Portugal'),
    @('RU', 'Russia'),
    @('SE', 'Sweden'),
    @('CH', 'Switzerland'),
    @('SK', 'Slovakia'),
    @('UA', 'Ukraine'),
    @('BE', 'Belgium'),
    @('AL', 'Albania'),
    @('GR', 'Greece'),
    @('XK', 'Kosovo')
)

for ($i = 0; $i -lt $sheet002.Length; $i++) {
    $r = $i + 1
    $row = $sheet002[$i]
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
}

# ---------------------------------------------------------------------
# 4. View state: zoom + selection per sheet, and which tab is active.
#    "000" keeps the filter view but is no longer the active tab;
#    "002" becomes the active (selected) tab, matching the diff.
# ---------------------------------------------------------------------
$ws0.Activate()
$excel.ActiveWindow.Zoom = 310
$ws0.Range("A1").Select()

$ws1.Activate()
$excel.ActiveWindow.Zoom = 310
$ws1.Range("B21").Select()

$ws2.Activate()
$excel.ActiveWindow.Zoom = 310
$ws2.Range("A30").Select()

Write-Output "Added sheets 001 and 002 with reference data."
